# Update the "Generate Report for Handback" timestamps.
# These cells are text values (formatted to look like dates) that record
# when the handoff/handback xliff files were generated.

$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for the first row.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-10-21 04:23:03"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K)
# for the first data row.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-10-21 04:22:51"
$wsZhCn.Range("K2").Value = "2016-10-21 04:23:34"

# de-de sheet: "Correspond Handoff Datetime" (H) for the first data row.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-10-21 04:23:52"
